$wb = $excel.ActiveWorkbook

# --- Update sunrise/sunset timestamp strings (shared across both sheets) ---
$daily = $wb.Worksheets.Item("Daily")
$hourly = $wb.Worksheets.Item("Hourly")

$daily.Range("E2").Value = "2024-02-06T07:42:17"
$daily.Range("F2").Value = "2024-02-06T17:37:38"

$hourly.Range("E2").Value = "2024-02-06T07:42:17"
$hourly.Range("F2").Value = "2024-02-06T17:37:38"
$hourly.Range("E3").Value = "2024-02-06T07:42:17"
$hourly.Range("F3").Value = "2024-02-06T17:37:38"
$hourly.Range("E4").Value = "2024-02-06T07:42:17"
$hourly.Range("F4").Value = "2024-02-06T17:37:38"
$hourly.Range("E5").Value = "2024-02-06T07:42:17"
$hourly.Range("F5").Value = "2024-02-06T17:37:38"
$hourly.Range("E6").Value = "2024-02-06T07:42:17"
$hourly.Range("F6").Value = "2024-02-06T17:37:38"
$hourly.Range("E7").Value = "2024-02-06T07:42:17"
$hourly.Range("F7").Value = "2024-02-06T17:37:38"
$hourly.Range("E8").Value = "2024-02-06T07:42:17"
$hourly.Range("F8").Value = "2024-02-06T17:37:38"
$hourly.Range("E9").Value = "2024-02-06T07:42:17"
$hourly.Range("F9").Value = "2024-02-06T17:37:38"
$hourly.Range("E10").Value = "2024-02-06T07:42:17"
$hourly.Range("F10").Value = "2024-02-06T17:37:38"
$hourly.Range("E11").Value = "2024-02-06T07:42:17"
$hourly.Range("F11").Value = "2024-02-06T17:37:38"
$hourly.Range("E12").Value = "2024-02-06T07:42:17"
$hourly.Range("F12").Value = "2024-02-06T17:37:38"
$hourly.Range("E13").Value = "2024-02-06T07:42:17"
$hourly.Range("F13").Value = "2024-02-06T17:37:38"
$hourly.Range("E14").Value = "2024-02-06T07:42:17"
$hourly.Range("F14").Value = "2024-02-06T17:37:38"
$hourly.Range("E15").Value = "2024-02-06T07:42:17"
$hourly.Range("F15").Value = "2024-02-06T17:37:38"
$hourly.Range("E16").Value = "2024-02-06T07:42:17"
$hourly.Range("F16").Value = "2024-02-06T17:37:38"
$hourly.Range("E17").Value = "2024-02-06T07:42:17"
$hourly.Range("F17").Value = "2024-02-06T17:37:38"
$hourly.Range("E18").Value = "2024-02-06T07:42:17"
$hourly.Range("F18").Value = "2024-02-06T17:37:38"
$hourly.Range("E19").Value = "2024-02-06T07:42:17"
$hourly.Range("F19").Value = "2024-02-06T17:37:38"
$hourly.Range("E20").Value = "2024-02-06T07:42:17"
$hourly.Range("F20").Value = "2024-02-06T17:37:38"
$hourly.Range("E21").Value = "2024-02-06T07:42:17"
$hourly.Range("F21").Value = "2024-02-06T17:37:38"
$hourly.Range("E22").Value = "2024-02-06T07:42:17"
$hourly.Range("F22").Value = "2024-02-06T17:37:38"
$hourly.Range("E23").Value = "2024-02-06T07:42:17"
$hourly.Range("F23").Value = "2024-02-06T17:37:38"
$hourly.Range("E24").Value = "2024-02-06T07:42:17"
$hourly.Range("F24").Value = "2024-02-06T17:37:38"
$hourly.Range("E25").Value = "2024-02-06T07:42:17"
$hourly.Range("F25").Value = "2024-02-06T17:37:38"

# --- Update Daily sheet row 2 (lat/lon + GHI/DNI/DHI values) ---
$daily.Range("A2").Value = 46.073272
$daily.Range("B2").Value = 23.580489
$daily.Range("G2").Value = 2707.88
$daily.Range("H2").Value = 5929.04
$daily.Range("I2").Value = 694.1900000000001
$daily.Range("J2").Value = 2542.34
$daily.Range("K2").Value = 4705.11
$daily.Range("L2").Value = 738.09

# --- Update Hourly sheet lat/lon for every row, and GHI/DNI/DHI for daylight hours ---
# row 2
$hourly.Range("A2").Value = 46.073272
$hourly.Range("B2").Value = 23.580489
# row 3
$hourly.Range("A3").Value = 46.073272
$hourly.Range("B3").Value = 23.580489
# row 4
$hourly.Range("A4").Value = 46.073272
$hourly.Range("B4").Value = 23.580489
# row 5
$hourly.Range("A5").Value = 46.073272
$hourly.Range("B5").Value = 23.580489
# row 6
$hourly.Range("A6").Value = 46.073272
$hourly.Range("B6").Value = 23.580489
# row 7
$hourly.Range("A7").Value = 46.073272
$hourly.Range("B7").Value = 23.580489
# row 8
$hourly.Range("A8").Value = 46.073272
$hourly.Range("B8").Value = 23.580489
# row 9
$hourly.Range("A9").Value = 46.073272
$hourly.Range("B9").Value = 23.580489
$hourly.Range("H9").Value = 1.16
$hourly.Range("I9").Value = 10.81
$hourly.Range("J9").Value = 2.74
$hourly.Range("K9").Value = 1.16
$hourly.Range("L9").Value = 0
$hourly.Range("M9").Value = 1.16
# row 10
$hourly.Range("A10").Value = 46.073272
$hourly.Range("B10").Value = 23.580489
$hourly.Range("H10").Value = 76.01000000000001
$hourly.Range("I10").Value = 352.91
$hourly.Range("J10").Value = 40.93
$hourly.Range("K10").Value = 76
$hourly.Range("L10").Value = 264.93
$hourly.Range("M10").Value = 39.64
# row 11
$hourly.Range("A11").Value = 46.073272
$hourly.Range("B11").Value = 23.580489
$hourly.Range("H11").Value = 214.47
$hourly.Range("I11").Value = 599.88
$hourly.Range("J11").Value = 67.70999999999999
$hourly.Range("K11").Value = 214.42
$hourly.Range("L11").Value = 593.86
$hourly.Range("M11").Value = 60.66
# row 12
$hourly.Range("A12").Value = 46.073272
$hourly.Range("B12").Value = 23.580489
$hourly.Range("H12").Value = 337.15
$hourly.Range("I12").Value = 711.35
$hourly.Range("J12").Value = 82.61
$hourly.Range("K12").Value = 336.92
$hourly.Range("L12").Value = 711.52
$hourly.Range("M12").Value = 74.81999999999999
# row 13
$hourly.Range("A13").Value = 46.073272
$hourly.Range("B13").Value = 23.580489
$hourly.Range("H13").Value = 421.56
$hourly.Range("I13").Value = 766.0700000000001
$hourly.Range("J13").Value = 90.83
$hourly.Range("K13").Value = 420.16
$hourly.Range("L13").Value = 760.09
$hourly.Range("M13").Value = 85
# row 14
$hourly.Range("A14").Value = 46.073272
$hourly.Range("B14").Value = 23.580489
$hourly.Range("H14").Value = 456.88
$hourly.Range("I14").Value = 785.65
$hourly.Range("J14").Value = 93.94
$hourly.Range("K14").Value = 451.96
$hourly.Range("L14").Value = 763.83
$hourly.Range("M14").Value = 92.26000000000001
# row 15
$hourly.Range("A15").Value = 46.073272
$hourly.Range("B15").Value = 23.580489
$hourly.Range("H15").Value = 439.04
$hourly.Range("I15").Value = 776.04
$hourly.Range("J15").Value = 92.38
$hourly.Range("K15").Value = 431.82
$hourly.Range("L15").Value = 745.9299999999999
$hourly.Range("M15").Value = 91.81999999999999
# row 16
$hourly.Range("A16").Value = 46.073272
$hourly.Range("B16").Value = 23.580489
$hourly.Range("H16").Value = 370.07
$hourly.Range("I16").Value = 734.34
$hourly.Range("J16").Value = 85.94
$hourly.Range("K16").Value = 346.31
$hourly.Range("L16").Value = 606.77
$hourly.Range("M16").Value = 103.34
# row 17
$hourly.Range("A17").Value = 46.073272
$hourly.Range("B17").Value = 23.580489
$hourly.Range("H17").Value = 258.39
$hourly.Range("I17").Value = 645.9400000000001
$hourly.Range("J17").Value = 73.55
$hourly.Range("K17").Value = 185.99
$hourly.Range("L17").Value = 213.56
$hourly.Range("M17").Value = 120.26
# row 18
$hourly.Range("A18").Value = 46.073272
$hourly.Range("B18").Value = 23.580489
$hourly.Range("H18").Value = 121.25
$hourly.Range("I18").Value = 461.61
$hourly.Range("J18").Value = 51.88
$hourly.Range("K18").Value = 71.28
$hourly.Range("L18").Value = 44.62
$hourly.Range("M18").Value = 62.81
# row 19
$hourly.Range("A19").Value = 46.073272
$hourly.Range("B19").Value = 23.580489
$hourly.Range("H19").Value = 11.9
$hourly.Range("I19").Value = 84.43000000000001
$hourly.Range("J19").Value = 11.68
$hourly.Range("K19").Value = 6.33
$hourly.Range("L19").Value = 0
$hourly.Range("M19").Value = 6.33
# row 20
$hourly.Range("A20").Value = 46.073272
$hourly.Range("B20").Value = 23.580489
# row 21
$hourly.Range("A21").Value = 46.073272
$hourly.Range("B21").Value = 23.580489
# row 22
$hourly.Range("A22").Value = 46.073272
$hourly.Range("B22").Value = 23.580489
# row 23
$hourly.Range("A23").Value = 46.073272
$hourly.Range("B23").Value = 23.580489
# row 24
$hourly.Range("A24").Value = 46.073272
$hourly.Range("B24").Value = 23.580489
# row 25
$hourly.Range("A25").Value = 46.073272
$hourly.Range("B25").Value = 23.580489
